$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.008355500328646
$ws.Cells.Item(2, 4).Value = 1.031862482820368
$ws.Cells.Item(2, 5).Value = 1.010850847508915
$ws.Cells.Item(2, 6).Value = 1.02340218746966
$ws.Cells.Item(2, 9).Value = 1.031136902952628
$ws.Cells.Item(2, 10).Value = 1.013621809153451
$ws.Cells.Item(2, 11).Value = 1.034669374833291
$ws.Cells.Item(2, 12).Value = 1.013719791725444
$ws.Cells.Item(2, 13).Value = 1.026233731250371
$ws.Cells.Item(2, 14).Value = 1.008457817609619

# Row 3
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.010119072505176
$ws.Cells.Item(3, 4).Value = 1.032308744751241
$ws.Cells.Item(3, 5).Value = 1.012374099252135
$ws.Cells.Item(3, 6).Value = 1.02515244823389
$ws.Cells.Item(3, 9).Value = 1.031231007818702
$ws.Cells.Item(3, 10).Value = 1.015012910850549
$ws.Cells.Item(3, 11).Value = 1.034925260177525
$ws.Cells.Item(3, 12).Value = 1.015045196086176
$ws.Cells.Item(3, 13).Value = 1.027788293297456
$ws.Cells.Item(3, 14).Value = 1.008941661028133

# Row 4
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.011257240933905
$ws.Cells.Item(4, 4).Value = 1.032596757619517
$ws.Cells.Item(4, 5).Value = 1.013357282464112
$ws.Cells.Item(4, 6).Value = 1.026281169441578
$ws.Cells.Item(4, 9).Value = 1.031289914643131
$ws.Cells.Item(4, 10).Value = 1.015909912525705
$ws.Cells.Item(4, 11).Value = 1.035089359344523
$ws.Cells.Item(4, 12).Value = 1.01589988997523
$ws.Cells.Item(4, 13).Value = 1.028789963106094
$ws.Cells.Item(4, 14).Value = 1.009253016127193

# Row 5
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.011735030463957
$ws.Cells.Item(5, 4).Value = 1.032717658114054
$ws.Cells.Item(5, 5).Value = 1.013770036635452
$ws.Cells.Item(5, 6).Value = 1.02675478593049
$ws.Cells.Item(5, 9).Value = 1.031314203978149
$ws.Cells.Item(5, 10).Value = 1.016286275059398
$ws.Cells.Item(5, 11).Value = 1.035157992971316
$ws.Cells.Item(5, 12).Value = 1.016258513960801
$ws.Cells.Item(5, 13).Value = 1.029210064905027
$ws.Cells.Item(5, 14).Value = 1.009383501561096

# Row 6
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.011815213080607
$ws.Cells.Item(6, 4).Value = 1.032737947249239
$ws.Cells.Item(6, 5).Value = 1.013839306486317
$ws.Cells.Item(6, 6).Value = 1.026834256010785
$ws.Cells.Item(6, 9).Value = 1.031318254409935
$ws.Cells.Item(6, 10).Value = 1.016349425153905
$ws.Cells.Item(6, 11).Value = 1.035169496109278
$ws.Cells.Item(6, 12).Value = 1.016318688442926
$ws.Cells.Item(6, 13).Value = 1.029280543594021
$ws.Cells.Item(6, 14).Value = 1.009405386828473

# Row 7
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.011263627894369
$ws.Cells.Item(7, 4).Value = 1.032598373806598
$ws.Cells.Item(7, 5).Value = 1.013362799950615
$ws.Cells.Item(7, 6).Value = 1.02628750143548
$ws.Cells.Item(7, 9).Value = 1.031290241064828
$ws.Cells.Item(7, 10).Value = 1.015914944377866
$ws.Cells.Item(7, 11).Value = 1.035090277820369
$ws.Cells.Item(7, 12).Value = 1.015904684618048
$ws.Cells.Item(7, 13).Value = 1.028795580437932
$ws.Cells.Item(7, 14).Value = 1.009254761277954

# Row 8
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.008952135320147
$ws.Cells.Item(8, 4).Value = 1.032013453544404
$ws.Cells.Item(8, 5).Value = 1.011366155602676
$ws.Cells.Item(8, 6).Value = 1.023994495374074
$ws.Cells.Item(8, 9).Value = 1.031169117328695
$ws.Cells.Item(8, 10).Value = 1.014092595504263
$ws.Cells.Item(8, 11).Value = 1.034756157838758
$ws.Cells.Item(8, 12).Value = 1.014168332821032
$ws.Cells.Item(8, 13).Value = 1.026759988211323
$ws.Cells.Item(8, 14).Value = 1.008621694415907

# Row 9
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.004855389828109
$ws.Cells.Item(9, 4).Value = 1.030977056932842
$ws.Cells.Item(9, 5).Value = 1.00782834472848
$ws.Cells.Item(9, 6).Value = 1.019924019591568
$ws.Cells.Item(9, 9).Value = 1.030940466761065
$ws.Cells.Item(9, 10).Value = 1.01085677761882
$ws.Cells.Item(9, 11).Value = 1.034156107733131
$ws.Cells.Item(9, 12).Value = 1.011085646497541
$ws.Cells.Item(9, 13).Value = 1.023139946812617
$ws.Cells.Item(9, 14).Value = 1.007492744386687

# Row 10
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.002107274366038
$ws.Cells.Item(10, 4).Value = 1.03028234931679
$ws.Cells.Item(10, 5).Value = 1.005455858928028
$ws.Cells.Item(10, 6).Value = 1.017189288261364
$ws.Cells.Item(10, 9).Value = 1.030777786946658
$ws.Cells.Item(10, 10).Value = 1.00868217432303
$ws.Cells.Item(10, 11).Value = 1.033748501543117
$ws.Cells.Item(10, 12).Value = 1.009014259468418
$ws.Cells.Item(10, 13).Value = 1.020703469636454
$ws.Cells.Item(10, 14).Value = 1.006730813217755

# Row 11
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.000913046507045
$ws.Cells.Item(11, 4).Value = 1.029980647278354
$ws.Cells.Item(11, 5).Value = 1.004425046928264
$ws.Cells.Item(11, 6).Value = 1.015999900643611
$ws.Cells.Item(11, 9).Value = 1.030704910146481
$ws.Cells.Item(11, 10).Value = 1.00773623209878
$ws.Cells.Item(11, 11).Value = 1.033570209960555
$ws.Cells.Item(11, 12).Value = 1.00811329144943
$ws.Cells.Item(11, 13).Value = 1.019642766682969
$ws.Cells.Item(11, 14).Value = 1.006398619203517

# Row 12
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.000468793701433
$ws.Cells.Item(12, 4).Value = 1.029868448838169
$ws.Cells.Item(12, 5).Value = 1.004041613481
$ws.Cells.Item(12, 6).Value = 1.015557304230258
$ws.Cells.Item(12, 9).Value = 1.030677474103791
$ws.Cells.Item(12, 10).Value = 1.007384200825833
$ws.Cells.Item(12, 11).Value = 1.033503714892501
$ws.Cells.Item(12, 12).Value = 1.007778008899381
$ws.Cells.Item(12, 13).Value = 1.019247901894614
$ws.Cells.Item(12, 14).Value = 1.006274880506629

# Row 13
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.000564117797152
$ws.Cells.Item(13, 4).Value = 1.029892521786199
$ws.Cells.Item(13, 5).Value = 1.004123886142333
$ws.Cells.Item(13, 6).Value = 1.015652279394395
$ws.Cells.Item(13, 9).Value = 1.030683375814145
$ws.Cells.Item(13, 10).Value = 1.007459743143274
$ws.Cells.Item(13, 11).Value = 1.033517990510538
$ws.Cells.Item(13, 12).Value = 1.007849956581637
$ws.Cells.Item(13, 13).Value = 1.019332641493822
$ws.Cells.Item(13, 14).Value = 1.006301438672378

# Row 14
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.00087633808741
$ws.Cells.Item(14, 4).Value = 1.029971375626867
$ws.Cells.Item(14, 5).Value = 1.004393363393806
$ws.Cells.Item(14, 6).Value = 1.015963332016516
$ws.Cells.Item(14, 9).Value = 1.030702649750944
$ws.Cells.Item(14, 10).Value = 1.007707146770599
$ws.Cells.Item(14, 11).Value = 1.033564718957481
$ws.Cells.Item(14, 12).Value = 1.008085589686167
$ws.Cells.Item(14, 13).Value = 1.019610144935641
$ws.Cells.Item(14, 14).Value = 1.006388398036531

# Row 15
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.001068618771538
$ws.Cells.Item(15, 4).Value = 1.030019942494427
$ws.Cells.Item(15, 5).Value = 1.004559324639151
$ws.Cells.Item(15, 6).Value = 1.016154874573577
$ws.Cells.Item(15, 9).Value = 1.030714476499617
$ws.Cells.Item(15, 10).Value = 1.007859491534319
$ws.Cells.Item(15, 11).Value = 1.033593474169123
$ws.Cells.Item(15, 12).Value = 1.008230688013255
$ws.Cells.Item(15, 13).Value = 1.019781007831156
$ws.Cells.Item(15, 14).Value = 1.006441930414943

# Row 16
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.002186439462589
$ws.Cells.Item(16, 4).Value = 1.030302353616663
$ws.Cells.Item(16, 5).Value = 1.005524195158826
$ws.Cells.Item(16, 6).Value = 1.01726811217439
$ws.Cells.Item(16, 9).Value = 1.030782572190095
$ws.Cells.Item(16, 10).Value = 1.008744860875578
$ws.Cells.Item(16, 11).Value = 1.033760296350582
$ws.Cells.Item(16, 12).Value = 1.009073967242618
$ws.Cells.Item(16, 13).Value = 1.020773743598681
$ws.Cells.Item(16, 14).Value = 1.006752811480127

# Row 17
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.002886459354645
$ws.Cells.Item(17, 4).Value = 1.03047926513213
$ws.Cells.Item(17, 5).Value = 1.006128481540855
$ws.Cells.Item(17, 6).Value = 1.017965002906894
$ws.Cells.Item(17, 9).Value = 1.030824634477808
$ws.Cells.Item(17, 10).Value = 1.009299060270662
$ws.Cells.Item(17, 11).Value = 1.03386445889757
$ws.Cells.Item(17, 12).Value = 1.009601840595528
$ws.Cells.Item(17, 13).Value = 1.021394924408792
$ws.Cells.Item(17, 14).Value = 1.006947206547493

# Row 18
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.003294358001694
$ws.Cells.Item(18, 4).Value = 1.030582368800631
$ws.Cells.Item(18, 5).Value = 1.006480614034928
$ws.Cells.Item(18, 6).Value = 1.018370984071417
$ws.Cells.Item(18, 9).Value = 1.030848933754344
$ws.Cells.Item(18, 10).Value = 1.00962189925981
$ws.Cells.Item(18, 11).Value = 1.03392504191299
$ws.Cells.Item(18, 12).Value = 1.009909351075547
$ws.Cells.Item(18, 13).Value = 1.021756700625343
$ws.Cells.Item(18, 14).Value = 1.007060374991582

# Row 19
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.003433371799731
$ws.Cells.Item(19, 4).Value = 1.030617509907539
$ws.Cells.Item(19, 5).Value = 1.006600625364901
$ws.Cells.Item(19, 6).Value = 1.018509328341289
$ws.Cells.Item(19, 9).Value = 1.03085717934269
$ws.Cells.Item(19, 10).Value = 1.009731908922896
$ws.Cells.Item(19, 11).Value = 1.033945669767917
$ws.Cells.Item(19, 12).Value = 1.010014138664067
$ws.Cells.Item(19, 13).Value = 1.02187996452819
$ws.Cells.Item(19, 14).Value = 1.007098925543445

# Row 20
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.002811396520054
$ws.Cells.Item(20, 4).Value = 1.030460293073078
$ws.Cells.Item(20, 5).Value = 1.006063682363147
$ws.Cells.Item(20, 6).Value = 1.01789028533368
$ws.Cells.Item(20, 9).Value = 1.030820145897175
$ws.Cells.Item(20, 10).Value = 1.009239643082806
$ws.Cells.Item(20, 11).Value = 1.033853301159925
$ws.Cells.Item(20, 12).Value = 1.009545245139502
$ws.Cells.Item(20, 13).Value = 1.021328334397712
$ws.Cells.Item(20, 14).Value = 1.006926372479666

# Row 21
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.000784415477227
$ws.Cells.Item(21, 4).Value = 1.029948158798569
$ws.Cells.Item(21, 5).Value = 1.004314024210789
$ws.Cells.Item(21, 6).Value = 1.01587175714504
$ws.Cells.Item(21, 9).Value = 1.030696984173993
$ws.Cells.Item(21, 10).Value = 1.007634311092881
$ws.Cells.Item(21, 11).Value = 1.033550966037106
$ws.Cells.Item(21, 12).Value = 1.008016218929819
$ws.Cells.Item(21, 13).Value = 1.01952845131342
$ws.Cells.Item(21, 14).Value = 1.006362800296703

# Row 22
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 0.9995061231764496
$ws.Cells.Item(22, 4).Value = 1.029625391231303
$ws.Cells.Item(22, 5).Value = 1.00321078902746
$ws.Cells.Item(22, 6).Value = 1.014597961883781
$ws.Cells.Item(22, 9).Value = 1.030617427648238
$ws.Cells.Item(22, 10).Value = 1.006621112334999
$ws.Cells.Item(22, 11).Value = 1.033359315945887
$ws.Cells.Item(22, 12).Value = 1.007051247723971
$ws.Cells.Item(22, 13).Value = 1.018391736725866
$ws.Cells.Item(22, 14).Value = 1.006006449354488

# Row 23
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.000184142027801
$ws.Cells.Item(23, 4).Value = 1.029796569037974
$ws.Cells.Item(23, 5).Value = 1.003795939601577
$ws.Cells.Item(23, 6).Value = 1.015273673593257
$ws.Cells.Item(23, 9).Value = 1.030659803168534
$ws.Cells.Item(23, 10).Value = 1.007158599875655
$ws.Cells.Item(23, 11).Value = 1.033461061127888
$ws.Cells.Item(23, 12).Value = 1.007563144756059
$ws.Cells.Item(23, 13).Value = 1.018994815780752
$ws.Cells.Item(23, 14).Value = 1.006195550212007

# Row 24
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.002845315453818
$ws.Cells.Item(24, 4).Value = 1.030468865995459
$ws.Cells.Item(24, 5).Value = 1.006092963364963
$ws.Cells.Item(24, 6).Value = 1.017924048544063
$ws.Cells.Item(24, 9).Value = 1.030822174819623
$ws.Cells.Item(24, 10).Value = 1.009266492437948
$ws.Cells.Item(24, 11).Value = 1.033858343396796
$ws.Cells.Item(24, 12).Value = 1.009570819392019
$ws.Cells.Item(24, 13).Value = 1.021358425250559
$ws.Cells.Item(24, 14).Value = 1.00693578717447

# Row 25
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.00591740433733
$ws.Cells.Item(25, 4).Value = 1.031245659532406
$ws.Cells.Item(25, 5).Value = 1.008745348216443
$ws.Cells.Item(25, 6).Value = 1.020979978651345
$ws.Cells.Item(25, 9).Value = 1.031001383079161
$ws.Cells.Item(25, 10).Value = 1.011696316697362
$ws.Cells.Item(25, 11).Value = 1.034312570939695
$ws.Cells.Item(25, 12).Value = 1.011885402701242
$ws.Cells.Item(25, 13).Value = 1.024079823218827
$ws.Cells.Item(25, 14).Value = 1.007786222720407
